$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 7 (shifts existing rows 7-31 down to 8-32)
$ws.Rows(7).Insert()

# Populate the newly inserted row 7 with the new weekly data point
$ws.Range("A7").Value2 = 7
$ws.Range("B7").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C7").Value2 = "Ñuble"
$ws.Range("D7").Value2 = 44635
$ws.Range("E7").Value2 = 16
$ws.Range("F7").Value2 = 100112040
$ws.Range("G7").Value2 = "Cilantro"
$ws.Range("H7").Value2 = "Sin especificar"
$ws.Range("I7").Value2 = "Primera"
$ws.Range("J7").Value2 = 120
$ws.Range("K7").Value2 = 550
$ws.Range("L7").Value2 = 600
$ws.Range("M7").Value2 = 575
$ws.Range("N7").Value2 = "$/atado 0,5 a 1 kilo"
$ws.Range("O7").Value2 = "Provincia de Diguillín"
$ws.Range("P7").Value2 = 575
$ws.Range("Q7").Value2 = 1
$ws.Range("R7").Value2 = "Hortaliza"
